# Revised 7-Day Training Program Agenda — reformat:
#  - strip bold ("Focus/Target Audience/Objective/Program Structure/Total
#    Duration/Delivery Format/Daily Session" labels are no longer bold)
#  - merge the label + value runs on the Focus/Target Audience/Objective/
#    Program Structure lines into single runs
#  - drop the bulleted-list formatting from "Total Duration", "Delivery
#    Format" and "Daily Session"
#  - replace the horizontal-rule picture paragraph with the "Delivery
#    Format" text, move "Daily Session" into its own paragraph right
#    after it, and add one more blank paragraph at the end

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Work from the tail of the document forward so paragraph indices for
#    the parts we still need to touch stay valid.
# ---------------------------------------------------------------------

# The old horizontal-rule paragraph (contains just a <w:pict> hr image).
$hrPara = $d.Paragraphs.Item(6)

# New paragraph for "Daily Session: ..." right after the hr paragraph.
$hrPara.Range.InsertParagraphAfter() | Out-Null
$dailyPara = $d.Paragraphs.Item(7)
$dailyRange = $dailyPara.Range
$dailyRange.MoveEnd(1, -1)
$dailyRange.Font.Bold = $false
$dailyRange.Font.BoldBi = $false
$dailyRange.Text = "Daily Session: 8 hours/day (inclusive of breaks)"

# One extra trailing blank paragraph (document already ends with one
# empty paragraph; the edit adds a second one before it).
$d.Paragraphs.Item(7).Range.InsertParagraphAfter() | Out-Null

# Turn the old hr paragraph into the "Delivery Format" line.
$hrRange = $hrPara.Range
$hrRange.MoveEnd(1, -1)
$hrRange.Font.Bold = $false
$hrRange.Font.BoldBi = $false
$hrRange.Text = "Delivery Format: 60% Hands-on, 20% Group Activities, 20% Theory"

# Drop the old numbered "Daily Session" paragraph (#5) -- its text now
# lives in the paragraph created above.
$d.Paragraphs.Item(5).Range.Delete() | Out-Null

# Drop the old numbered "Delivery Format" paragraph (#4) -- its text now
# lives in the paragraph that used to hold the hr picture.
$d.Paragraphs.Item(4).Range.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2) "Total Duration" (#3): remove the bullet numbering and the bold
#    label formatting, merge into a single run.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.ListFormat.RemoveNumbers()
$p3Range = $p3.Range
$p3Range.MoveEnd(1, -1)
$p3Range.Font.Bold = $false
$p3Range.Font.BoldBi = $false
$p3Range.Text = "Total Duration: 7 Days"

# ---------------------------------------------------------------------
# 3) Second paragraph: merge "Focus:"/"Target Audience:"/"Objective:"/
#    "Program Structure:" labels with their values into single runs,
#    keep the line breaks, remove all bold formatting.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2Range = $p2.Range
$p2Range.MoveEnd(1, -1)
$p2Range.Font.Bold = $false
$p2Range.Font.BoldBi = $false
$lineBreak = [string][char]11
$p2Text = "Focus: SDLC, C#, Golang, Angular, and Testing" + $lineBreak + `
    "Target Audience: 25 Freshers or Beginners" + $lineBreak + `
    "Objective: Equip participants with foundational skills in software development and full-stack web development." + $lineBreak + `
    "Program Structure:"
$p2Range.Text = $p2Text

# ---------------------------------------------------------------------
# 4) Title paragraph (#1): remove bold.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Range = $p1.Range
$p1Range.MoveEnd(1, -1)
$p1Range.Font.Bold = $false
$p1Range.Font.BoldBi = $false
$p1Range.Text = "Revised 7-Day Training Program Agenda"
